$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values ---
$ws.Range("A2").Value = 76683052018
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 249.99
$ws.Range("A3").Value = 741569060080
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = 487.99

# --- Apply grey Arial font style to the UPC cells (A2:A3) ---
# Build the format on a scratch cell first so the whole A2:A3 range
# picks up a single, shared style entry instead of one-per-cell.
$ws.Range("E1").Font.Name = "Arial"
$ws.Range("E1").Font.Color = 7697781
$ws.Range("E1").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$ws.Range("E1").Clear()

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 29.917
$ws.Columns("B").ColumnWidth = 32.917
$ws.Columns("C").ColumnWidth = 17.917

# --- Selection ---
$ws.Range("D3").Select()
